$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.91
$ws.Range("I2").Value = 4.2
$ws.Range("J2").Value = 2.63
$ws.Range("L2").Value = 4.75
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.65
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.5
$ws.Range("U2").Value = 1.95
$ws.Range("V2").Value = 1.8
$ws.Range("W2").Value = 6.5
$ws.Range("X2").Value = 8.5
$ws.Range("Z2").Value = 17
$ws.Range("AF2").Value = 51
$ws.Range("AH2").Value = 10
$ws.Range("AK2").Value = 41
$ws.Range("AO2").Value = 11
$ws.Range("AQ2").Value = 41
$ws.Range("AR2").Value = 67
$ws.Range("AT2").Value = 2.5
$ws.Range("AU2").Value = 8.5
$ws.Range("AX2").Value = 23

# Row 4 updates
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
